$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sending / target cluster names in the fixed order used by the source data
$clusters = @("ECs", "FAPs", "M2", "sCs")
$targets  = @("ECs", "FAPs", "sCs")

# Per sending-cluster aggregate stats (columns E-J), indexed by sending cluster
$sendStats = @{
    "ECs"  = @{ E=3; F=1; G=98.47161033333335;  H=295.414831;  I=0.272681344498213;  J=0.2726813444982129 }
    "FAPs" = @{ E=3; F=1; G=113.1680936666667;  H=339.504281;  I=0.3133779150241075; J=0.3133779150241075 }
    "M2"   = @{ E=3; F=1; G=58.24795766666667;  H=174.743873;  I=0.1612965539718111; J=0.1612965539718111 }
    "sCs"  = @{ E=3; F=1; G=91.235724;          H=273.707172;  I=0.2526441865058685; J=0.2526441865058685 }
}

# Per target-cluster aggregate stats (columns K-P), indexed by target cluster
$targetStats = @{
    "ECs"  = @{ K=3; L=1;                  M=43.55927533333334; N=130.677826; O=0.9894183625413969;   P=0.9894183625413967 }
    "FAPs" = @{ K=2; L=0.6666666666666666; M=0.205596;          N=0.616788;   O=0.004669968820840217; P=0.004669968820840216 }
    "sCs"  = @{ K=2; L=0.6666666666666666; M=0.260262;          N=0.780786;   O=0.005911668637762975; P=0.005911668637762974 }
}

# Per (sending,target) edge weight/specificity values (columns Q-T)
$edgeStats = @{
    "ECs|ECs"   = @{ Q=4289.35198702638;   R=38604.16788323742;  S=0.2697959293690084;    T=0.2697959293690083 }
    "ECs|FAPs"  = @{ Q=20.245369198092;    R=182.208322782828;   S=0.001273413376831445;  T=0.001273413376831444 }
    "ECs|sCs"   = @{ Q=25.628418248574;    R=230.655764237166;   S=0.001612001752373127;  T=0.001612001752373127 }
    "FAPs|ECs"  = @{ Q=4929.52015097479;   R=44365.68135877311;  S=0.3100618635397895;    T=0.3100618635397894 }
    "FAPs|FAPs" = @{ Q=23.266907385492;    R=209.402166469428;   S=0.001463465092302497;  T=0.001463465092302497 }
    "FAPs|sCs"  = @{ Q=29.453354393874;    R=265.080189544866;   S=0.001852586392015567;  T=0.001852586392015567 }
    "M2|ECs"    = @{ Q=2537.238825606678;  R=22835.1494304601;   S=0.1595897723143594;    T=0.1595897723143594 }
    "M2|FAPs"   = @{ Q=11.975547104436;    R=107.779923939924;   S=0.0007532498779573293; T=0.0007532498779573291 }
    "M2|sCs"    = @{ Q=15.159729958242;    R=136.437569624178;   S=0.0009535317794943989; T=0.0009535317794943987 }
    "sCs|ECs"   = @{ Q=3974.162021952009;  R=35767.45819756808;  S=0.2499707973182397;    T=0.2499707973182396 }
    "sCs|FAPs"  = @{ Q=18.757699911504;    R=168.819299203536;   S=0.001179840473748947;  T=0.001179840473748946 }
    "sCs|sCs"   = @{ Q=23.745191999688;    R=213.706727997192;   S=0.001493548713879882;  T=0.001493548713879882 }
}

$row = 2
foreach ($send in $clusters) {
    foreach ($targ in $targets) {
        $ws.Cells.Item($row, 1).Value = $send
        $ws.Cells.Item($row, 2).Value = "Anxa2"
        $ws.Cells.Item($row, 3).Value = "Robo4"
        $ws.Cells.Item($row, 4).Value = $targ

        $ss = $sendStats[$send]
        $ws.Cells.Item($row, 5).Value = $ss.E
        $ws.Cells.Item($row, 6).Value = $ss.F
        $ws.Cells.Item($row, 7).Value = $ss.G
        $ws.Cells.Item($row, 8).Value = $ss.H
        $ws.Cells.Item($row, 9).Value = $ss.I
        $ws.Cells.Item($row, 10).Value = $ss.J

        $ts = $targetStats[$targ]
        $ws.Cells.Item($row, 11).Value = $ts.K
        $ws.Cells.Item($row, 12).Value = $ts.L
        $ws.Cells.Item($row, 13).Value = $ts.M
        $ws.Cells.Item($row, 14).Value = $ts.N
        $ws.Cells.Item($row, 15).Value = $ts.O
        $ws.Cells.Item($row, 16).Value = $ts.P

        $es = $edgeStats["$send|$targ"]
        $ws.Cells.Item($row, 17).Value = $es.Q
        $ws.Cells.Item($row, 18).Value = $es.R
        $ws.Cells.Item($row, 19).Value = $es.S
        $ws.Cells.Item($row, 20).Value = $es.T

        $row = $row + 1
    }
}
